$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the two discontinued products (rows 56 and 57 in the original sheet):
# "Фп Фиточай "Баланс" (БАД) 20х2,0 г" and "Фп Фиточай "Дивный вечер" (БАД) 20х2,0 г"
$ws.Rows("56:57").Delete()

# Refresh the stock figures (and re-sequence a few rows) for the rest of the table
$data = @(
    @('ВердиоГаст® Растительный комплекс для улучшения пищеварения (БАД ),  капсулы', 81536),
    @('Пижма цветки 75г', 11501),
    @('Полынь горькая трава 50г', 26124),
    @('Алтей корни 75г', 3888),
    @('Чистотел трава 50г', 10513),
    @('Сенна листья 50г', 14489),
    @('Шиповник плоды низковитаминные 50г', 22321),
    @('Липа цветки 35г', 14594),
    @('Кукуруза столбики с рыльцами 40г', 20974),
    @('Брусника листья 50г', 13185),
    @('Дуба кора 75г', 56245),
    @('Эвкалипт прутовидный листья 75г', 23651),
    @('Мята перечная листья 50г', 21023),
    @('Укроп пахучий плоды 50г', 54488),
    @('Сб. Фитонефрол (Урологический сбор) 50г', 7521),
    @('Ноготки цветки 50г', 23015),
    @('Багульник болотный побеги 50г', 13305),
    @('Эрва шерстистая трава 30г', 12668),
    @('Ромашка цветки вн 50г', 91497),
    @('Сб. Грудной №4 50г', 35220),
    @('Девясил корневища и корни 50г', 18671),
    @('Валериана корневища с корнями 50г', 21306),
    @('Береза почки 50г', 18963),
    @('Пустырник трава 50г', 12278),
    @('Боярышник плоды 75г', 23160),
    @('Бессмертник песчаный цветки 30г', 30775),
    @('Чабрец трава 50г', 22009),
    @('Подорожник большой листья 50г', 10178),
    @('Шалфей листья 50г', 40146),
    @('Сб. Фитогепатол №2 (Желчегонный сбор №2) 35г', 5027),
    @('Аир корневища 75г', 8569),
    @('Мать-и-мачеха листья 35г', 32930),
    @('Рябина плоды 50г', 2156),
    @('Ламинарии слоевища (морская капуста) 100г', 19008),
    @('Лен семена 100г', 68099),
    @('Крушина кора 50г', 12676),
    @('Череда трава 50г', 14983),
    @('Зверобой трава 50г', 47815),
    @('Солодка корни 50г', 52318),
    @('Тысячелистник трава 50г', 21641),
    @('Можжевельник плоды 50г', 18772),
    @('Толокнянка листья 50г', 11409),
    @('Сб. Фитопектол №1 (Грудной сбор №1) 35г', 8693),
    @('Спорыш трава 50г', 23725),
    @('Чага (березовый гриб) 50г', 47698),
    @('Крапива листья 50г', 24295),
    @('Сб. Фитопектол №2 (Грудной сбор №2) 35г', 11334),
    @('Фп Детский травяной чай "ФармаЦветик® для иммунитета" 20х1,5 г', 1598),
    @('Фп Фиточай "Лактафитол" (БАД) 20х1,5 г', 11129),
    @('Фп Детский травяной чай "ФармаЦветик®  при простуде" 20х1,5 г', 3741),
    @('Фп Детский травяной чай "ФармаЦветик® для животика" 20х1,5 г', 4070),
    @('Фп Детский травяной чай "ФармаЦветик® для спокойного сна" 20х1,5 г', 6398),
    @('Фп "ВердиоГаст® Фиточай для улучшения пищеварения с зеленым чаем"(БАД) 20*1,5г', 8110),
    @('Фп "ВердиоГаст® Фиточай для улучшения пищеварения с черным чаем" (БАД) 20*1,5г', 9670),
    @('Фп "Щедрость природы® Фиточай для иммунитета" 20х2,0 г', 1008),
    @('Фп "Щедрость природы® Фиточай при простуде" 20х2,0 г', 1026),
    @('Фп "Щедрость природы® Фиточай кардиологический" 20х2,0 г', 1476),
    @('Фп "Щедрость природы® Фиточай успокоительный"20х2,0 г', 1998),
    @('Фп Шалфей листья 20х1,5г', 101778),
    @('Фп Сб. Фитогепатол №3 (Желчегонный сбор №3) 20x2,0г', 49883),
    @('Фп "Щедрость природы® Фиточай диабетический" 20х2,0 г', 1062),
    @('Фп Сб. Фитоседан №2 (Успокоительный сбор №2) 20x2,0г', 40041),
    @('Фп Сб. Грудной №4 20x2,0г', 561621),
    @('Фп Сб. Бруснивер 20x2,0г', 155217),
    @('Фп Сб. Проктофитол (Противогеморроидальный сбор) 20х2,0г', 18262),
    @('Фп Толокнянка листья 20x1,5г', 28674),
    @('Фп Сб. Фитогастрол (Желудочно-кишечный сбор) 20x2,0г', 73431),
    @('Фп Липа цветки 20x1,5г', 64024),
    @('Фп Боярышник плоды 20х3,0г', 16352),
    @('Фп Аир корневища 20x1,5г', 3460),
    @('Фп Мелисса лекарственная трава 20x1,5г', 36937),
    @('Фп Сенна листья 20x1,5г', 66606),
    @('Фп Ромашка цветки 20x1,5г', 1365725),
    @('Фп Сб. Элекасол 20x2,0г', 42512),
    @('Фп Сб. Фитоседан №3 (Успокоительный сбор №3) 20х2,0г', 85659),
    @('Фп Пижма цветки 20х1,5г', 6006),
    @('Фп Фиточай "Тибетский" (БАД) 20х2,0  г', 7632),
    @('Фп Фиточай "Опалиховский" (БАД) 20х2,0 г', 5130),
    @('Фп "Щедрость природы® Фиточай очищающий" 20х2,0 г', 1818),
    @('Фп Сб. Арфазетин-Э 20x2,0г', 45377),
    @('Фп Подорожник листья 20x1,5г', 33337),
    @('Фп Шиповник плоды 20х2,0г', 57809),
    @('Фп Чистотел трава 20х1,5г', 34962),
    @('Фп Зверобой трава 20x1,5г', 61909),
    @('Фп Мята перечная листья 20x1,5г', 77626),
    @('Фп Брусника листья 20х1,5г', 89135),
    @('Фп Пустырник трава 20x1,5г', 49389),
    @('Фп "Щедрость природы® Фиточай для пищеварения" 20х2,0 г', 1854),
    @('Фп Череда трава 20х1,5г', 60192),
    @('Фп Сб. Фитонефрол (Урологический сбор) 20x2,0г', 229603),
    @('Фп Душица трава 20x1,5г', 36234),
    @('Фп Чабрец трава 20x1,5 г', 89835),
    @('Фп Сб. Желудочный №3 20x2,0г', 31917),
    @('Фп Хвощ полевой трава 20х1,5г', 37186),
    @('Фп Крапива листья 20x1,5г', 88216),
    @('Фп Пастушья сумка трава 20х1,5г', 6928),
    @('Фп Береза листья 20x1,5г', 6664),
    @('Фп Золототысячник трава 20х1,5г', 6685),
    @('Фп Фиалка трехцветная трава 20x1,5г', 6286),
    @('Фп Ольха соплодия 20х1,5г', 6181),
    @('Фп Почечный чай листья 20x1,5г', 47022),
    @('Фп Валериана корневища с корнями 20x1,5г', 16971),
    @('Фп Лапчатка корневища 20x2,5г', 3646),
    @('Фп Тысячелистник трава 20x1,5г', 18772),
    @('Фп Крушина кора 20x1,5г', 9711),
    @('Фп Девясил корневища и корни 20х1,5г', 14154),
    @('Фп Ноготки цветки 20x1,5г', 46164),
    @('Фп Бадан корневища 20x1,5г', 2311),
    @('Фп Кровохлебка корневища и корни 20x1,5г', 12454),
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $rowNum = $i + 2
    $ws.Cells.Item($rowNum, 1).Value = $data[$i][0]
    $ws.Cells.Item($rowNum, 2).Value = $data[$i][1]
}

# Restore the scroll position / active-cell selection recorded for the refreshed sheet
$ws.Application.ActiveWindow.ScrollRow = 76
[void]$ws.Range("A93").Select()
